# Insert a new data row at row 141 (shifts existing rows 141-162 down to 142-163)
# and populate it with the new observation, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(141).Insert()

$ws.Cells.Item(141, 1).Value = 4
$ws.Cells.Item(141, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(141, 3).Value = "Los Lagos"
$ws.Cells.Item(141, 4).Value = 44491
$ws.Cells.Item(141, 5).Value = 10
$ws.Cells.Item(141, 6).Value = 100112043
$ws.Cells.Item(141, 7).Value = "Pepino ensalada"
$ws.Cells.Item(141, 8).Value = "Sin especificar"
$ws.Cells.Item(141, 9).Value = "Primera"
$ws.Cells.Item(141, 10).Value = 400
$ws.Cells.Item(141, 11).Value = 13000
$ws.Cells.Item(141, 12).Value = 13000
$ws.Cells.Item(141, 13).Value = 13000
$ws.Cells.Item(141, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(141, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(141, 16).Value = 217
$ws.Cells.Item(141, 17).Value = 60
$ws.Cells.Item(141, 18).Value = "Hortaliza"
